$wb = $excel.ActiveWorkbook

# --- AddCustomerTest (sheet1): add new data rows 3-5 and a Runmode column (E) ---
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

$ws1.Range("A3").Value = "Virat "
$ws1.Range("B3").Value = "Kohli"
$ws1.Range("C3").Value = 12345
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "Rakesh"
$ws1.Range("B4").Value = "Roshan"
$ws1.Range("C4").Value = 879012
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "Arjun"
$ws1.Range("B5").Value = "Kapoor"
$ws1.Range("C5").Value = 231843
$ws1.Range("D5").Value = "Customer added successfully"

# --- test_suite (renamed from Sheet3): build the run-mode control table ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "test_suite"

$ws3.Range("A1").Value = "TCID"
$ws3.Range("A2").Value = "BankManagerLoginTest"
$ws3.Range("A3").Value = "OpenAccountTest"
$ws3.Range("A4").Value = "AddCustomerTest"

$ws3.Range("B2").Value = "Y"
$ws3.Range("B1").Value = "Runmode"
$ws3.Range("B3").Value = "N"
$ws3.Range("B4").Value = "Y"

$ws3.Columns.Item(1).ColumnWidth = 20.92

# --- Back on AddCustomerTest: add the Runmode header/values (column E) ---
$ws1.Range("E1").Value = "runmode"
$ws1.Range("E2").Value = "Y"
$ws1.Range("E3").Value = "N"
$ws1.Range("E4").Value = "Y"
$ws1.Range("E5").Value = "Y"

# --- Selections: test_suite gets a resting selection on B4, then AddCustomerTest
#     becomes the active tab with E1 selected, and OpenAccountTest stops being
#     the active tab. ---
[void]$ws3.Range("B4").Select()
[void]$ws1.Activate()
[void]$ws1.Range("E1").Select()
